# Auto-generated edit script: updates crypto price/volume data (and reorders
# a couple of rows) to match the scraped coinranking.com snapshot referenced
# in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'94.116.29"
$ws.Range("E2").Value = "'  +2.58%  "
$ws.Range("D3").Value = "'3.072.12"
$ws.Range("E3").Value = "'  -0.45%  "
$ws.Range("E4").Value = "'  +0.09%  "
$ws.Range("D5").Value = "'236.15"
$ws.Range("E5").Value = "'  +0.52%  "
$ws.Range("D6").Value = "'609.35"
$ws.Range("E6").Value = "'  -0.07%  "
$ws.Range("E7").Value = "'  +2.15%  "
$ws.Range("D8").Value = "'0.377"
$ws.Range("E8").Value = "'  -1.62%  "
$ws.Range("E9").Value = "'  +0.07%  "
$ws.Range("D10").Value = "'0.809"
$ws.Range("E10").Value = "'  +11.50%  "
$ws.Range("D11").Value = "'3.074.51"
$ws.Range("E11").Value = "'  -0.29%  "
$ws.Range("E12").Value = "'  -1.75%  "
$ws.Range("D13").Value = "'94.029.58"
$ws.Range("E13").Value = "'  +2.15%  "
$ws.Range("D14").Value = "'0.0000240"
$ws.Range("E14").Value = "'  -2.32%  "
$ws.Range("D15").Value = "'33.92"
$ws.Range("E15").Value = "'  +0.37%  "
$ws.Range("D16").Value = "'5.31"
$ws.Range("E16").Value = "'  -1.38%  "
$ws.Range("D17").Value = "'3.656.53"
$ws.Range("E17").Value = "'  -0.39%  "
$ws.Range("D18").Value = "'3.070.88"
$ws.Range("D19").Value = "'3.55"
$ws.Range("E19").Value = "'  -2.81%  "
$ws.Range("D20").Value = "'14.36"
$ws.Range("E20").Value = "'  -1.16%  "
$ws.Range("D21").Value = "'5.74"
$ws.Range("E21").Value = "'  -0.11%  "
$ws.Range("D22").Value = "'444.42"
$ws.Range("E22").Value = "'  +0.68%  "
$ws.Range("D23").Value = "'8.82"
$ws.Range("E23").Value = "'  -4.47%  "
$ws.Range("D24").Value = "'0.0000189"
$ws.Range("E24").Value = "'  -1.95%  "
$ws.Range("D25").Value = "'8.35"
$ws.Range("E25").Value = "'  +6.70%  "
$ws.Range("D26").Value = "'5.51"
$ws.Range("E26").Value = "'  -3.19%  "
$ws.Range("B27").Value = "'Litecoin"
$ws.Range("C27").Value = "'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D27").Value = "'84.74"
$ws.Range("E27").Value = "'  -1.04%  "
$ws.Range("B28").Value = "'Aptos"
$ws.Range("C28").Value = "'https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'11.93"
$ws.Range("E28").Value = "'  +3.60%  "
$ws.Range("D29").Value = "'3.242.92"
$ws.Range("E29").Value = "'  -0.45%  "
$ws.Range("E30").Value = "'  +0.25%  "
$ws.Range("D31").Value = "'0.251"
$ws.Range("E31").Value = "'  +9.60%  "
$ws.Range("E32").Value = "'  +6.99%  "
$ws.Range("D33").Value = "'0.122"
$ws.Range("E33").Value = "'  -5.93%  "
$ws.Range("B34").Value = "'Binance-PegBSC-USD"
$ws.Range("C34").Value = "'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").Value = "'0.999"
$ws.Range("E34").Value = "'  +0.67%  "
$ws.Range("B35").Value = "'InternetComputer(DFINITY)"
$ws.Range("C35").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D35").Value = "'8.94"
$ws.Range("E35").Value = "'  -1.15%  "
$ws.Range("D36").Value = "'7.58"
$ws.Range("E36").Value = "'  -2.64%  "
$ws.Range("D37").Value = "'25.43"
$ws.Range("E37").Value = "'  -1.39%  "
$ws.Range("D38").Value = "'0.150"
$ws.Range("E38").Value = "'  -4.38%  "
$ws.Range("D39").Value = "'1.88"
$ws.Range("E39").Value = "'  +0.01%  "
$ws.Range("D40").Value = "'481.31"
$ws.Range("E40").Value = "'  +0.20%  "
$ws.Range("E41").Value = "'  +0.76%  "
$ws.Range("D42").Value = "'0.436"
$ws.Range("E42").Value = "'  +1.93%  "
$ws.Range("E43").Value = "'  -4.75%  "
$ws.Range("D44").Value = "'1.24"
$ws.Range("E44").Value = "'  -2.76%  "
$ws.Range("D46").Value = "'3.07"
$ws.Range("E46").Value = "'  -6.17%  "
$ws.Range("D47").Value = "'161.49"
$ws.Range("E47").Value = "'  +0.55%  "
$ws.Range("D48").Value = "'0.671"
$ws.Range("E48").Value = "'  -1.35%  "
$ws.Range("D49").Value = "'1.81"
$ws.Range("E49").Value = "'  -2.50%  "
$ws.Range("D50").Value = "'43.63"
$ws.Range("E50").Value = "'  -0.41%  "
$ws.Range("D51").Value = "'0.998"
$ws.Range("E51").Value = "'  +0.16%  "
